$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.734.39"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "3.145.24"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'569.73"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").Value = "'149.32"
$ws.Range("E6").Value = "  +5.40%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.143.08"
$ws.Range("E8").Value = "  +3.41%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  +5.38%  "
$ws.Range("E10").Value = "  +5.97%  "
$ws.Range("D11").Value = "'6.19"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "'0.504"
$ws.Range("E12").Value = "  +8.11%  "
$ws.Range("E13").Value = "  +15.84%  "
$ws.Range("D14").Value = "'38.28"
$ws.Range("E14").Value = "  +10.83%  "
$ws.Range("D15").Value = "3.656.30"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "64.788.22"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").Value = "'7.20"
$ws.Range("E17").Value = "  +8.30%  "
$ws.Range("D18").Value = "3.143.84"
$ws.Range("E18").Value = "  +3.10%  "
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'514.01"
$ws.Range("E20").Value = "  +8.36%  "
$ws.Range("D21").Value = "'14.89"
$ws.Range("E21").Value = "  +7.63%  "
$ws.Range("D22").Value = "'0.736"
$ws.Range("E22").Value = "  +10.13%  "
$ws.Range("D23").Value = "'15.41"
$ws.Range("E23").Value = "  +10.45%  "
$ws.Range("E24").Value = "  +5.23%  "
$ws.Range("D25").Value = "'84.91"
$ws.Range("E25").Value = "  +5.32%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'2.92"
$ws.Range("E27").Value = "  +5.76%  "
$ws.Range("E28").Value = "  +12.70%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  +7.96%  "
$ws.Range("D30").Value = "'27.83"
$ws.Range("E30").Value = "  +7.17%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'2.69"
$ws.Range("E32").Value = "  +10.41%  "
$ws.Range("E33").Value = "  +4.68%  "
$ws.Range("E34").Value = "  +11.17%  "
$ws.Range("D35").Value = "'6.59"
$ws.Range("E35").Value = "  +7.67%  "
$ws.Range("D36").Value = "'55.65"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").Value = "'485.30"
$ws.Range("E37").Value = "  +12.00%  "
$ws.Range("D38").Value = "'0.0863"
$ws.Range("E38").Value = "  +8.00%  "
$ws.Range("D39").Value = "'0.0424"
$ws.Range("E39").Value = "  +5.25%  "
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("D41").Value = "3.113.32"
$ws.Range("E41").Value = "  +5.90%  "
$ws.Range("D42").Value = "'8.65"
$ws.Range("E42").Value = "  +6.74%  "
$ws.Range("D43").Value = "'0.120"
$ws.Range("E43").Value = "  +6.58%  "
$ws.Range("E44").Value = "  +14.27%  "
$ws.Range("E45").Value = "  +17.12%  "
$ws.Range("D46").Value = "'29.50"
$ws.Range("E46").Value = "  +6.19%  "
$ws.Range("D47").Value = "0.0₃0573"
$ws.Range("E47").Value = "  +13.05%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'0.116"
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("E50").Value = "  +12.10%  "
$ws.Range("D51").Value = "'120.21"
$ws.Range("E51").Value = "  +2.79%  "
